{"js": "// Blank out the document: remove every paragraph except one, leaving a\n// single empty paragraph (mirrors a \"select all, delete\" in the UI).\nconst body = context.document.body;\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nwhile (paragraphs.items.length > 1) {\n  paragraphs.items[0].delete();\n  await context.sync();\n\n  paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n}\n\n// Deleting the paragraphs also removes the hidden \"_GoBack\" bookmark\n// Word leaves at the site of the last edit, so restore it on the\n// now-empty remaining paragraph.\nconst start = body.getRange(\"Start\");\nstart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Collapse the whole body down to a single empty paragraph, the way\n# \"select all, delete\" does: repeatedly drop the first paragraph's\n# content until only one paragraph mark is left.\nwhile ($d.Paragraphs.Count -gt 1) {\n  $d.Content.Delete()\n}\n\n# The bulk delete above also wipes out the hidden \"_GoBack\" bookmark\n# that Word drops at the site of the last edit, so put it back at the\n# (now empty) remaining paragraph.\n$r = $d.Paragraphs(1).Range\n$r.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
